$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values for rows 2-19 as part of a data repull /
# mean recalculation. Row 17's value (0) is unchanged by this repull.
$values = @{
    2  = 3
    3  = -3
    4  = -4
    5  = 6
    6  = -1
    7  = 1
    8  = -6
    9  = -1
    10 = -3
    11 = -4
    12 = -2
    13 = 1
    14 = 2
    15 = 4
    16 = 2
    18 = -3
    19 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
